# Update TPM-derived ligand/receptor expression & specificity values for Igf1-Igf1r
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.447974666666667
$ws.Range("H2").Value = 28.343924
$ws.Range("I2").Value = 0.06074879557565148
$ws.Range("J2").Value = 0.06110214141073961
$ws.Range("M2").Value = 8.021311666666668
$ws.Range("N2").Value = 24.063935
$ws.Range("O2").Value = 0.2318437811880539
$ws.Range("P2").Value = 0.2575223477274881
$ws.Range("Q2").Value = 75.78514942010446
$ws.Range("R2").Value = 682.0663447809401
$ws.Range("S2").Value = 0.01408423046887916
$ws.Range("T2").Value = 0.01573516690727064
$ws.Range("G3").Value = 9.447974666666667
$ws.Range("H3").Value = 28.343924
$ws.Range("I3").Value = 0.06074879557565148
$ws.Range("J3").Value = 0.06110214141073961
$ws.Range("O3").Value = 0.3382003894878138
$ws.Range("P3").Value = 0.3756588072233373
$ws.Range("Q3").Value = 110.5510224166071
$ws.Range("R3").Value = 994.959201749464
$ws.Range("S3").Value = 0.02054526632460091
$ws.Range("T3").Value = 0.02295355756115012
$ws.Range("G4").Value = 9.447974666666667
$ws.Range("H4").Value = 28.343924
$ws.Range("I4").Value = 0.06074879557565148
$ws.Range("J4").Value = 0.06110214141073961
$ws.Range("M4").Value = 2.475956
$ws.Range("N4").Value = 7.427868
$ws.Range("O4").Value = 0.07156373233578578
$ws.Range("P4").Value = 0.07948999222154987
$ws.Range("Q4").Value = 23.39276956378134
$ws.Range("R4").Value = 210.534926074032
$ws.Range("S4").Value = 0.00434741054629729
$ws.Range("T4").Value = 0.004857008745459732
$ws.Range("G5").Value = 9.447974666666667
$ws.Range("H5").Value = 28.343924
$ws.Range("I5").Value = 0.06074879557565148
$ws.Range("J5").Value = 0.06110214141073961
$ws.Range("M5").Value = 10.349683
$ws.Range("N5").Value = 20.699366
$ws.Range("O5").Value = 0.2991418038011306
$ws.Range("P5").Value = 0.2215161123395049
$ws.Range("Q5").Value = 97.78354279203066
$ws.Range("R5").Value = 586.701256752184
$ws.Range("S5").Value = 0.01817250428724653
$ws.Range("T5").Value = 0.01353510882092571
$ws.Range("G6").Value = 9.447974666666667
$ws.Range("H6").Value = 28.343924
$ws.Range("I6").Value = 0.06074879557565148
$ws.Range("J6").Value = 0.06110214141073961
$ws.Range("M6").Value = 2.049936666666667
$ws.Range("N6").Value = 6.14981
$ws.Range("O6").Value = 0.05925029318721588
$ws.Range("P6").Value = 0.06581274048811983
$ws.Range("Q6").Value = 19.36774969493778
$ws.Range("R6").Value = 174.30974725444
$ws.Range("S6").Value = 0.003599383948627593
$ws.Range("T6").Value = 0.004021299375933406
$ws.Range("I7").Value = 0.4644313471136943
$ws.Range("J7").Value = 0.4671327156039162
$ws.Range("M7").Value = 8.021311666666668
$ws.Range("N7").Value = 24.063935
$ws.Range("O7").Value = 0.2318437811880539
$ws.Range("P7").Value = 0.2575223477274881
$ws.Range("Q7").Value = 579.3859565916878
$ws.Range("R7").Value = 5214.47360932519
$ws.Range("S7").Value = 0.1076755196171005
$ws.Range("T7").Value = 0.1202971136226375
$ws.Range("I8").Value = 0.4644313471136943
$ws.Range("J8").Value = 0.4671327156039162
$ws.Range("O8").Value = 0.3382003894878138
$ws.Range("P8").Value = 0.3756588072233373
$ws.Range("S8").Value = 0.1570708624842015
$ws.Range("T8").Value = 0.1754825187587656
$ws.Range("I9").Value = 0.4644313471136943
$ws.Range("J9").Value = 0.4671327156039162
$ws.Range("M9").Value = 2.475956
$ws.Range("N9").Value = 7.427868
$ws.Range("O9").Value = 0.07156373233578578
$ws.Range("P9").Value = 0.07948999222154987
$ws.Range("Q9").Value = 178.8403437183813
$ws.Range("R9").Value = 1609.563093465432
$ws.Range("S9").Value = 0.03323644061319284
$ws.Range("T9").Value = 0.03713237592978677
$ws.Range("I10").Value = 0.4644313471136943
$ws.Range("J10").Value = 0.4671327156039162
$ws.Range("M10").Value = 10.349683
$ws.Range("N10").Value = 20.699366
$ws.Range("O10").Value = 0.2991418038011306
$ws.Range("P10").Value = 0.2215161123395049
$ws.Range("Q10").Value = 747.5661381285805
$ws.Range("R10").Value = 4485.396828771483
$ws.Range("S10").Value = 0.1389308309173795
$ws.Range("T10").Value = 0.1034774231071751
$ws.Range("I11").Value = 0.4644313471136943
$ws.Range("J11").Value = 0.4671327156039162
$ws.Range("M11").Value = 2.049936666666667
$ws.Range("N11").Value = 6.14981
$ws.Range("O11").Value = 0.05925029318721588
$ws.Range("P11").Value = 0.06581274048811983
$ws.Range("Q11").Value = 148.0686159477711
$ws.Range("R11").Value = 1332.61754352994
$ws.Range("S11").Value = 0.02751769348182001
$ws.Range("T11").Value = 0.03074328418555122
$ws.Range("G12").Value = 20.46218966666666
$ws.Range("H12").Value = 61.38656899999999
$ws.Range("I12").Value = 0.131568237738417
$ws.Range("J12").Value = 0.1323335054016559
$ws.Range("M12").Value = 8.021311666666668
$ws.Range("N12").Value = 24.063935
$ws.Range("O12").Value = 0.2318437811880539
$ws.Range("P12").Value = 0.2575223477274881
$ws.Range("Q12").Value = 164.1336006987794
$ws.Range("R12").Value = 1477.202406289015
$ws.Range("S12").Value = 0.0305032777215234
$ws.Range("T12").Value = 0.03407883499404266
$ws.Range("G13").Value = 20.46218966666666
$ws.Range("H13").Value = 61.38656899999999
$ws.Range("I13").Value = 0.131568237738417
$ws.Range("J13").Value = 0.1323335054016559
$ws.Range("O13").Value = 0.3382003894878138
$ws.Range("P13").Value = 0.3756588072233373
$ws.Range("Q13").Value = 239.428667872437
$ws.Range("R13").Value = 2154.858010851934
$ws.Range("S13").Value = 0.04449642924735792
$ws.Range("T13").Value = 0.04971224679486912
$ws.Range("G14").Value = 20.46218966666666
$ws.Range("H14").Value = 61.38656899999999
$ws.Range("I14").Value = 0.131568237738417
$ws.Range("J14").Value = 0.1323335054016559
$ws.Range("M14").Value = 2.475956
$ws.Range("N14").Value = 7.427868
$ws.Range("O14").Value = 0.07156373233578578
$ws.Range("P14").Value = 0.07948999222154987
$ws.Range("Q14").Value = 50.66348127832133
$ws.Range("R14").Value = 455.971331504892
$ws.Range("S14").Value = 0.009415514149403103
$ws.Range("T14").Value = 0.01051918931502805
$ws.Range("G15").Value = 20.46218966666666
$ws.Range("H15").Value = 61.38656899999999
$ws.Range("I15").Value = 0.131568237738417
$ws.Range("J15").Value = 0.1323335054016559
$ws.Range("M15").Value = 10.349683
$ws.Range("N15").Value = 20.699366
$ws.Range("O15").Value = 0.2991418038011306
$ws.Range("P15").Value = 0.2215161123395049
$ws.Range("Q15").Value = 211.7771765358756
$ws.Range("R15").Value = 1270.663059215254
$ws.Range("S15").Value = 0.03935755996000605
$ws.Range("T15").Value = 0.02931400364883369
$ws.Range("G16").Value = 20.46218966666666
$ws.Range("H16").Value = 61.38656899999999
$ws.Range("I16").Value = 0.131568237738417
$ws.Range("J16").Value = 0.1323335054016559
$ws.Range("M16").Value = 2.049936666666667
$ws.Range("N16").Value = 6.14981
$ws.Range("O16").Value = 0.05925029318721588
$ws.Range("P16").Value = 0.06581274048811983
$ws.Range("Q16").Value = 41.94619287798778
$ws.Range("R16").Value = 377.51573590189
$ws.Range("S16").Value = 0.007795456660126527
$ws.Range("T16").Value = 0.008709230648882384
$ws.Range("G17").Value = 2.6981485
$ws.Range("H17").Value = 5.396297000000001
$ws.Range("I17").Value = 0.01734861464410334
$ws.Range("J17").Value = 0.01163301532943533
$ws.Range("M17").Value = 8.021311666666668
$ws.Range("N17").Value = 24.063935
$ws.Range("O17").Value = 0.2318437811880539
$ws.Range("P17").Value = 0.2575223477274881
$ws.Range("Q17").Value = 21.64269004144917
$ws.Range("R17").Value = 129.856140248695
$ws.Range("S17").Value = 0.004022168417463361
$ws.Range("T17").Value = 0.002995761418786044
$ws.Range("G18").Value = 2.6981485
$ws.Range("H18").Value = 5.396297000000001
$ws.Range("I18").Value = 0.01734861464410334
$ws.Range("J18").Value = 0.01163301532943533
$ws.Range("O18").Value = 0.3382003894878138
$ws.Range("P18").Value = 0.3756588072233373
$ws.Range("Q18").Value = 31.57111294542367
$ws.Range("R18").Value = 189.426677672542
$ws.Range("S18").Value = 0.00586730822970974
$ws.Range("T18").Value = 0.004370044663066474
$ws.Range("G19").Value = 2.6981485
$ws.Range("H19").Value = 5.396297000000001
$ws.Range("I19").Value = 0.01734861464410334
$ws.Range("J19").Value = 0.01163301532943533
$ws.Range("M19").Value = 2.475956
$ws.Range("N19").Value = 7.427868
$ws.Range("O19").Value = 0.07156373233578578
$ws.Range("P19").Value = 0.07948999222154987
$ws.Range("Q19").Value = 6.680496967466
$ws.Range("R19").Value = 40.08298180479601
$ws.Range("S19").Value = 0.001241531614787305
$ws.Range("T19").Value = 0.0009247082980499848
$ws.Range("G20").Value = 2.6981485
$ws.Range("H20").Value = 5.396297000000001
$ws.Range("I20").Value = 0.01734861464410334
$ws.Range("J20").Value = 0.01163301532943533
$ws.Range("M20").Value = 10.349683
$ws.Range("N20").Value = 20.699366
$ws.Range("O20").Value = 0.2991418038011306
$ws.Range("P20").Value = 0.2215161123395049
$ws.Range("Q20").Value = 27.9249816619255
$ws.Range("R20").Value = 111.699926647702
$ws.Range("S20").Value = 0.005189695878087782
$ws.Range("T20").Value = 0.002576900330562379
$ws.Range("G21").Value = 2.6981485
$ws.Range("H21").Value = 5.396297000000001
$ws.Range("I21").Value = 0.01734861464410334
$ws.Range("J21").Value = 0.01163301532943533
$ws.Range("M21").Value = 2.049936666666667
$ws.Range("N21").Value = 6.14981
$ws.Range("O21").Value = 0.05925029318721588
$ws.Range("P21").Value = 0.06581274048811983
$ws.Range("Q21").Value = 5.531033542261668
$ws.Range("R21").Value = 33.18620125357
$ws.Range("S21").Value = 0.00102791050405515
$ws.Range("T21").Value = 0.0007656006189704471
$ws.Range("G22").Value = 50.68616266666667
$ws.Range("H22").Value = 152.058488
$ws.Range("I22").Value = 0.3259030049281339
$ws.Range("J22").Value = 0.3277986222542529
$ws.Range("M22").Value = 8.021311666666668
$ws.Range("N22").Value = 24.063935
$ws.Range("O22").Value = 0.2318437811880539
$ws.Range("P22").Value = 0.2575223477274881
$ws.Range("Q22").Value = 406.5695079366978
$ws.Range("R22").Value = 3659.12557143028
$ws.Range("S22").Value = 0.07555858496308752
$ws.Range("T22").Value = 0.08441547078475123
$ws.Range("G23").Value = 50.68616266666667
$ws.Range("H23").Value = 152.058488
$ws.Range("I23").Value = 0.3259030049281339
$ws.Range("J23").Value = 0.3277986222542529
$ws.Range("O23").Value = 0.3382003894878138
$ws.Range("P23").Value = 0.3756588072233373
$ws.Range("Q23").Value = 593.0802423659965
$ws.Range("R23").Value = 5337.722181293968
$ws.Range("S23").Value = 0.1102205232019438
$ws.Range("T23").Value = 0.1231404394454859
$ws.Range("G24").Value = 50.68616266666667
$ws.Range("H24").Value = 152.058488
$ws.Range("I24").Value = 0.3259030049281339
$ws.Range("J24").Value = 0.3277986222542529
$ws.Range("M24").Value = 2.475956
$ws.Range("N24").Value = 7.427868
$ws.Range("O24").Value = 0.07156373233578578
$ws.Range("P24").Value = 0.07948999222154987
$ws.Range("Q24").Value = 125.4967085715093
$ws.Range("R24").Value = 1129.470377143584
$ws.Range("S24").Value = 0.02332283541210525
$ws.Range("T24").Value = 0.02605670993322532
$ws.Range("G25").Value = 50.68616266666667
$ws.Range("H25").Value = 152.058488
$ws.Range("I25").Value = 0.3259030049281339
$ws.Range("J25").Value = 0.3277986222542529
$ws.Range("M25").Value = 10.349683
$ws.Range("N25").Value = 20.699366
$ws.Range("O25").Value = 0.2991418038011306
$ws.Range("P25").Value = 0.2215161123395049
$ws.Range("Q25").Value = 524.5857160864347
$ws.Range("R25").Value = 3147.514296518608
$ws.Range("S25").Value = 0.09749121275841074
$ws.Range("T25").Value = 0.07261267643200801
$ws.Range("G26").Value = 50.68616266666667
$ws.Range("H26").Value = 152.058488
$ws.Range("I26").Value = 0.3259030049281339
$ws.Range("J26").Value = 0.3277986222542529
$ws.Range("M26").Value = 2.049936666666667
$ws.Range("N26").Value = 6.14981
$ws.Range("O26").Value = 0.05925029318721588
$ws.Range("P26").Value = 0.06581274048811983
$ws.Range("Q26").Value = 103.9034233430311
$ws.Range("R26").Value = 935.1308100872801
$ws.Range("S26").Value = 0.0193098485925866
$ws.Range("T26").Value = 0.02157332565878237
